$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.298.92'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '1.908.37'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.693'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.68%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.80'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.352'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0732'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0995'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '2.190.45'
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.711'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.35%  '
$ws.Range('D16').Value = '1.926.92'
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.86'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').Value = '35.335.67'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.17%  '
$ws.Range('D20').Value = '0.0₃0826'
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '241.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.58'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +21.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.20%  '
$ws.Range('E30').Value = '  +4.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.985'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.21'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0570'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.35'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.11'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0661'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +13.03%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0210'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.01%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '90.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').Value = '1.350.19'
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.66%  '
$ws.Range('B47').Value = 'Gas'
$ws.Range('C47').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.44%  '
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '47.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.80'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.82%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.41'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.88%  '

